$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").Value = "2024-09-25T18:06:40Z"
$ws.Range("B88").Value = "temperature"

# "25" must be stored as text (matching the rest of the sheet, which keeps
# every value as a string), not auto-coerced to a number by Excel.
$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "25"
$ws.Range("C88").ClearFormats()

$ws.Range("D88").Value = "N/A"
$ws.Range("E88").Value = "N/A"
$ws.Range("F88").Value = "N/A"
